# Set E10:E14 (in_service column for extr1..extr5) to TRUE
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E10").Value = $true
$ws.Range("E11").Value = $true
$ws.Range("E12").Value = $true
$ws.Range("E13").Value = $true
$ws.Range("E14").Value = $true
